$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# New handoff timestamp for the 39d2f9db-3fa7-4df7-ad84-00b724b5d7a6 row (row 6 in each sheet)
$zhcn.Range("E6").Value = "2016-03-22 18:40:19"
$dede.Range("E6").Value = "2016-03-22 18:40:23"
$overview.Range("D6").Value = "2016-03-22 18:40:23"
